# Update the division-problem worksheet numbers in the single table.
# Each populated row of the table (rows 1, 5, 9, 13, 17) holds five
# "aa÷b=" style expressions; replace them in place by cell so that
# duplicate original values (e.g. "18÷4=") are mapped independently.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "80÷3="
$t.Cell(1,2).Range.Text  = "15÷8="
$t.Cell(1,3).Range.Text  = "94÷2="
$t.Cell(1,4).Range.Text  = "53÷4="
$t.Cell(1,5).Range.Text  = "36÷8="

$t.Cell(5,1).Range.Text  = "69÷2="
$t.Cell(5,2).Range.Text  = "97÷5="
$t.Cell(5,3).Range.Text  = "34÷2="
$t.Cell(5,4).Range.Text  = "84÷6="
$t.Cell(5,5).Range.Text  = "75÷7="

$t.Cell(9,1).Range.Text  = "15÷9="
$t.Cell(9,2).Range.Text  = "85÷4="
$t.Cell(9,3).Range.Text  = "85÷7="
$t.Cell(9,4).Range.Text  = "63÷7="
$t.Cell(9,5).Range.Text  = "59÷3="

$t.Cell(13,1).Range.Text = "71÷3="
$t.Cell(13,2).Range.Text = "12÷2="
$t.Cell(13,3).Range.Text = "36÷5="
$t.Cell(13,4).Range.Text = "29÷5="
$t.Cell(13,5).Range.Text = "81÷2="

$t.Cell(17,1).Range.Text = "43÷8="
$t.Cell(17,2).Range.Text = "10÷7="
$t.Cell(17,3).Range.Text = "50÷9="
$t.Cell(17,4).Range.Text = "28÷8="
$t.Cell(17,5).Range.Text = "46÷7="
